# Add an "Exact?" column (H) to the color-analysis sheet:
#  - H2 header "Exact?" (bold + centered, matching the other header cells)
#  - H3:H15 "Yes"/"No" flags
#  - E6:G15 HSV->(OpenCV scale) formulas, mirroring rows 3:5 (filled down)
#  - selection left on B6, ready for more data entry

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell H2 (bold + centered, like the other header cells) --
$ws.Range("H2").Value = "Exact?"
$ws.Range("H2").HorizontalAlignment = -4108   # xlCenter
$ws.Range("H2").Font.Bold = $true

# --- Fill down the OpenCV-scale formulas for rows 6:15 --------------
$ws.Range("E6:E15").Formula = "=(B6/360) * 179"
$ws.Range("F6:F15").Formula = "=(C6/100) * 255"
$ws.Range("G6:G15").Formula = "=(D6/100) * 255"

# --- "Exact?" values for rows 3:15 -----------------------------------
$ws.Range("H3").Value = "No"
$ws.Range("H4").Value = "Yes"
$ws.Range("H5").Value = "No"
for ($r = 6; $r -le 15; $r++) {
    $ws.Cells.Item($r, 8).Value = "No"
}

$ws.Range("H3:H15").HorizontalAlignment = -4108   # xlCenter

# --- Leave selection on B6, as in the saved workbook -----------------
$ws.Range("B6").Select()
